$d = $word.ActiveDocument

# Locate the paragraph that ends the "Otras decisiones" bullet list with the
# admin/w23q user description - the new bullet item belongs right after it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Creamos el usuario con username: admin y password: w23q*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph to anchor the new bullet after."
}

# Insert a new paragraph right after the target one. InsertParagraphAfter
# clones the paragraph formatting (style + numbering) of $target, so the new
# paragraph automatically keeps the "Prrafodelista" style and numId=1 list.
$target.Range.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Range.Text = "Las publicaciones duran 1 mes"
